$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new paragraph right after the introductory paragraph ("... I
#    BILAGA 1 finns artfakta om fridlysta arter.") with the "Vi förväntar
#    oss..." text.
# ---------------------------------------------------------------------------
$introMatch = $d.Content
$introMatch.Find.Execute("I BILAGA 1 finns artfakta om fridlysta arter.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$introMatch.Collapse(0) | Out-Null
$insertionPoint = $introMatch.End
$introMatch.InsertParagraphAfter()

$newParaRange = $d.Range($insertionPoint + 1, $insertionPoint + 1)
$newParaRange.Text = "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden."

# ---------------------------------------------------------------------------
# 2. Remove the two empty paragraphs and the duplicate "Vi förväntar oss..."
#    paragraph that used to directly follow the second "... artskyddsförord-
#    ningen" comment paragraph (right before the page-break paragraph).
# ---------------------------------------------------------------------------
$searchText = "Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen"

$firstMatch = $d.Content
$firstMatch.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$secondMatch = $d.Range($firstMatch.End, $d.Content.End)
$secondMatch.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$commentPara = $d.Range($secondMatch.Start, $secondMatch.Start)
$commentPara.Expand(4) | Out-Null

$emptyPara1 = $d.Range($commentPara.End, $commentPara.End)
$emptyPara1.Expand(4) | Out-Null

$emptyPara2 = $d.Range($emptyPara1.End, $emptyPara1.End)
$emptyPara2.Expand(4) | Out-Null

$dupPara = $d.Range($emptyPara2.End, $emptyPara2.End)
$dupPara.Expand(4) | Out-Null

$d.Range($emptyPara1.Start, $dupPara.End).Delete()

# ---------------------------------------------------------------------------
# 3. Update the date shown in the "first page" header from 2023-11-13 to
#    2023-11-14.
# ---------------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$firstPageHeader = $sec.Headers.Item(2)
$headerDateRange = $firstPageHeader.Range
$headerDateRange.Find.Execute("2023-11-13", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headerDateRange.Text = "2023-11-14"
